$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in results (resultado/profit) for previously pending rows ---
function Set-Result($row, $resultado, $profit) {
    $ws.Cells.Item($row, 7).Value = $resultado
    $ws.Cells.Item($row, 8).Value = $profit
}

Set-Result 82 "Fallo"   -1
Set-Result 90 "Fallo"   -1
Set-Result 91 "Fallo"   -1
Set-Result 94 "Fallo"   -1
Set-Result 95 "Acierto" 1
Set-Result 96 "Fallo"   -1

# --- Append new match rows (106-111) ---
function Add-Row($row, $eventId, $fecha, $jugadorA, $jugadorB, $pronostico, $cuota) {
    # Seed the new row from the last existing row (105), which still has its
    # "resultado"/"profit" columns blank (pending match). This keeps G/H
    # as blank cells instead of the result columns being skipped outright.
    $ws.Range("A105:H105").Copy($ws.Range("A" + $row + ":H" + $row))

    $ws.Cells.Item($row, 1).Value = $eventId

    # "fecha" is stored as plain text (e.g. "2025-09-13"), not a date serial,
    # so force text formatting before assigning to avoid Excel's automatic
    # date recognition turning it into a date value.
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $fecha

    $ws.Cells.Item($row, 3).Value = $jugadorA
    $ws.Cells.Item($row, 4).Value = $jugadorB
    $ws.Cells.Item($row, 5).Value = $pronostico
    $ws.Cells.Item($row, 6).Value = $cuota
    # Columns 7/8 (resultado/profit) stay blank, as copied from row 105 —
    # these matches haven't been settled yet.
}

Add-Row 106 14673267 "2025-09-13" "Jessika Ponchet"      "Maria Carle"           "Gana Maria Carle"             2
Add-Row 107 14601387 "2025-09-13" "Stefano Napolitano"   "Gianluca Cadenasso"    "Gana Gianluca Cadenasso"      2.25
Add-Row 108 14601340 "2025-09-13" "Billy Harris"         "Juan Manuel Cerundolo" "Gana Juan Manuel Cerundolo"   2
Add-Row 109 14604884 "2025-09-13" "Mili Poljičak"        "Jay Clarke"            "Gana Jay Clarke"              1.8
Add-Row 110 14672757 "2025-09-13" "Akira Santillan"      "Kuan-Yi Lee"           "Gana Kuan-Yi Lee"             4.33
Add-Row 111 14674384 "2025-09-13" "Giles Hussey"         "Guillaume Dalmasso"    "Gana Guillaume Dalmasso"      5

$ws.Range("A1").Select() | Out-Null
